$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "51542051"

$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "51542054"

$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "51542055"

$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "06-21-2022"

$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "51542056"
